# Applies the "new input files generation" reshuffle to the bedrooms memory trial table.
# For each data row (2-42) this rewrites columns H (category) through V (r_perceptual)
# to the post-edit values captured in the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_f4jxo.png"; M=82.91666666666667; N=65.52777777777777; O=74.22222222222223; P=36.0; Q=8.0; R=8.0; S=8.0; T=8.0; U=8.0; V=8.0 },
    @{ Row=3; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_ce55l.png"; M=82.23809523809524; N=64.07142857142857; O=73.1547619047619; P=42.0; Q=8.0; R=8.0; S=8.0; T=8.0; U=8.0; V=8.0 },
    @{ Row=4; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_0nckg.png"; M=65.94285714285714; N=41.17142857142857; O=53.55714285714285; P=35.0; Q=4.0; R=4.0; S=4.0; T=4.0; U=4.0; V=4.0 },
    @{ Row=5; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_scrdm.png"; M=78.675; N=57.9; O=68.2875; P=40.0; Q=7.0; R=7.0; S=7.0; T=7.0; U=7.0; V=7.0 },
    @{ Row=6; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_bpyv5.png"; M=59.05882352941177; N=37.55882352941177; O=48.30882352941177; P=34.0; Q=3.0; R=3.0; S=3.0; T=3.0; U=3.0; V=3.0 },
    @{ Row=7; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_2js6m.png"; M=40.02777777777778; N=20.88888888888889; O=30.45833333333334; P=36.0; Q=2.0; R=2.0; S=2.0; T=2.0; U=2.0; V=2.0 },
    @{ Row=8; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_ozxpp.png"; M=26.26470588235294; N=11.47058823529412; O=18.86764705882353; P=34.0; Q=1.0; R=1.0; S=1.0; T=1.0; U=1.0; V=1.0 },
    @{ Row=9; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_uxxo0.png"; M=71.74418604651163; N=48.44186046511628; O=60.09302325581395; P=43.0; Q=5.0; R=5.0; S=5.0; T=5.0; U=5.0; V=5.0 },
    @{ Row=10; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_wijef.png"; M=69.875; N=48.025; O=58.95; P=40.0; Q=5.0; R=5.0; S=5.0; T=5.0; U=5.0; V=5.0 },
    @{ Row=11; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_a8y4y.png"; M=75.15789473684211; N=53.76315789473684; O=64.46052631578948; P=38.0; Q=6.0; R=6.0; S=6.0; T=6.0; U=6.0; V=6.0 },
    @{ Row=12; H=$null; I=$null; J="catch"; K="f"; L="stimuli/catch_20.jpg"; M=$null; N=$null; O=$null; P=$null; Q=$null; R=$null; S=$null; T=$null; U=$null; V=$null },
    @{ Row=13; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_cogrz.png"; M=60.5; N=39.71428571428572; O=50.10714285714286; P=42.0; Q=3.0; R=3.0; S=3.0; T=3.0; U=3.0; V=3.0 },
    @{ Row=14; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_a9acb.png"; M=77.11428571428571; N=58.42857142857143; O=67.77142857142857; P=35.0; Q=7.0; R=7.0; S=7.0; T=7.0; U=7.0; V=7.0 },
    @{ Row=15; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_oou46.png"; M=75.70270270270271; N=54.86486486486486; O=65.28378378378379; P=37.0; Q=6.0; R=6.0; S=6.0; T=6.0; U=6.0; V=6.0 },
    @{ Row=16; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_1vq1v.png"; M=69.42857142857143; N=46.59523809523809; O=58.01190476190476; P=42.0; Q=5.0; R=5.0; S=5.0; T=5.0; U=5.0; V=5.0 },
    @{ Row=17; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_cmyvx.png"; M=64.25; N=40.09375; O=52.171875; P=32.0; Q=4.0; R=4.0; S=4.0; T=4.0; U=4.0; V=4.0 },
    @{ Row=18; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_psgf7.png"; M=26.0; N=11.66666666666667; O=18.83333333333333; P=36.0; Q=1.0; R=1.0; S=1.0; T=1.0; U=1.0; V=1.0 },
    @{ Row=19; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_ca8kd.png"; M=92.05405405405405; N=73.02702702702703; O=82.54054054054055; P=37.0; Q=10.0; R=10.0; S=10.0; T=10.0; U=10.0; V=10.0 },
    @{ Row=20; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_e26ut.png"; M=81.07692307692308; N=61.28205128205128; O=71.17948717948718; P=39.0; Q=8.0; R=8.0; S=8.0; T=8.0; U=8.0; V=8.0 },
    @{ Row=21; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_juob3.png"; M=79.92105263157895; N=59.78947368421053; O=69.85526315789474; P=38.0; Q=7.0; R=7.0; S=7.0; T=7.0; U=7.0; V=7.0 },
    @{ Row=22; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_4wq98.png"; M=78.48387096774194; N=58.12903225806452; O=68.30645161290323; P=31.0; Q=7.0; R=7.0; S=7.0; T=7.0; U=7.0; V=7.0 },
    @{ Row=23; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_zt893.png"; M=68.53191489361703; N=49.19148936170212; O=58.86170212765958; P=47.0; Q=5.0; R=5.0; S=5.0; T=5.0; U=5.0; V=5.0 },
    @{ Row=24; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_kljj4.png"; M=64.35; N=44.15; O=54.25; P=40.0; Q=4.0; R=4.0; S=4.0; T=4.0; U=4.0; V=4.0 },
    @{ Row=25; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_3h4c9.png"; M=85.47619047619048; N=67.26190476190476; O=76.36904761904762; P=42.0; Q=9.0; R=9.0; S=9.0; T=9.0; U=9.0; V=9.0 },
    @{ Row=26; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_le8uf.png"; M=12.88888888888889; N=9.222222222222221; O=11.05555555555556; P=36.0; Q=1.0; R=1.0; S=1.0; T=1.0; U=1.0; V=1.0 },
    @{ Row=27; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_zi682.png"; M=84.6; N=69.525; O=77.0625; P=40.0; Q=9.0; R=9.0; S=9.0; T=9.0; U=9.0; V=9.0 },
    @{ Row=28; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_0eflx.png"; M=76.05128205128206; N=53.53846153846154; O=64.7948717948718; P=39.0; Q=6.0; R=6.0; S=6.0; T=6.0; U=6.0; V=6.0 },
    @{ Row=29; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_jp28n.png"; M=65.02564102564102; N=44.97435897435897; O=55.0; P=39.0; Q=4.0; R=4.0; S=4.0; T=4.0; U=4.0; V=5.0 },
    @{ Row=30; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_fnu4h.png"; M=85.87179487179488; N=70.71794871794872; O=78.2948717948718; P=39.0; Q=9.0; R=9.0; S=9.0; T=9.0; U=9.0; V=9.0 },
    @{ Row=31; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_gqy6z.png"; M=86.47368421052632; N=68.42105263157895; O=77.44736842105263; P=38.0; Q=9.0; R=9.0; S=9.0; T=9.0; U=9.0; V=9.0 },
    @{ Row=32; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_5yhyk.png"; M=46.375; N=31.325; O=38.85; P=40.0; Q=2.0; R=2.0; S=2.0; T=2.0; U=2.0; V=2.0 },
    @{ Row=33; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_g2akb.png"; M=87.875; N=79.0; O=83.4375; P=40.0; Q=10.0; R=10.0; S=10.0; T=10.0; U=10.0; V=10.0 },
    @{ Row=34; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_x0u5z.png"; M=92.0; N=78.16216216216216; O=85.08108108108108; P=37.0; Q=10.0; R=10.0; S=10.0; T=10.0; U=10.0; V=10.0 },
    @{ Row=35; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_72fmj.png"; M=53.87179487179487; N=36.02564102564103; O=44.94871794871795; P=39.0; Q=3.0; R=3.0; S=3.0; T=3.0; U=3.0; V=3.0 },
    @{ Row=36; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_t2ioc.png"; M=88.1891891891892; N=74.05405405405405; O=81.12162162162161; P=37.0; Q=10.0; R=10.0; S=10.0; T=10.0; U=10.0; V=10.0 },
    @{ Row=37; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_c4uwt.png"; M=44.48387096774194; N=30.06451612903226; O=37.2741935483871; P=31.0; Q=2.0; R=2.0; S=2.0; T=2.0; U=2.0; V=2.0 },
    @{ Row=38; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_zgg62.png"; M=82.1842105263158; N=63.52631578947368; O=72.85526315789474; P=38.0; Q=8.0; R=8.0; S=8.0; T=8.0; U=8.0; V=8.0 },
    @{ Row=39; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_5il0t.png"; M=48.09523809523809; N=30.90476190476191; O=39.5; P=42.0; Q=2.0; R=2.0; S=2.0; T=2.0; U=2.0; V=2.0 },
    @{ Row=40; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_rvssl.png"; M=74.25; N=54.33333333333334; O=64.29166666666667; P=36.0; Q=6.0; R=6.0; S=6.0; T=6.0; U=6.0; V=6.0 },
    @{ Row=41; H="bedrooms"; I="target"; J="old"; K="j"; L="stimuli/img_wyctg.png"; M=33.44736842105263; N=11.39473684210526; O=22.42105263157895; P=38.0; Q=1.0; R=1.0; S=1.0; T=1.0; U=1.0; V=1.0 },
    @{ Row=42; H="bedrooms"; I=$null; J="new"; K="f"; L="stimuli/img_u2o6z.png"; M=58.6; N=38.2; O=48.40000000000001; P=30.0; Q=3.0; R=3.0; S=3.0; T=3.0; U=3.0; V=3.0 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $cols = @{ 8="H"; 9="I"; 10="J"; 11="K"; 12="L"; 13="M"; 14="N"; 15="O"; 16="P"; 17="Q"; 18="R"; 19="S"; 20="T"; 21="U"; 22="V" }
    foreach ($colIndex in $cols.Keys) {
        $key = $cols[$colIndex]
        $value = $r[$key]
        if ($null -eq $value) {
            $ws.Cells.Item($row, $colIndex).ClearContents()
        } else {
            $ws.Cells.Item($row, $colIndex).Value = $value
        }
    }
}
